$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "Récupération de la liste du fichier JSON, avec leur image, le nom, le prix, et le descriptif"

# Row 3
$ws.Range("C3").Value = "Chaque produit est cliquable et effectue une redirection de page"

# Row 4
$ws.Range("C4").Value = "Récupération de l'image, le nom, le prix, et le descriptif"

# Row 5
$ws.Range("C5").Value = "Menu déroulant pour les coloris, et input numérique pour la quantité"
$ws.Range("E5").Value = "OK / Si aucun choix de couleur, ou quantité inférieure à 1 et supérieure à 100, popup alertant l'utilisateur d'indiquer sa sélection"

# Row 6
$ws.Range("C6").Value = "Envoie dans le Local Storage de la couleur et de la quantité"

# Row 7
$ws.Range("C7").Value = "Récupération des données du Local Storage, et récupération de l'image dans le fichier JSON"
$ws.Range("D7").Value = "Récapitulatif de tous les articles sélectionnés avec son coloris, sa quantité"

# Row 8
$ws.Range("C8").Value = "Input numérique ajoutant ou diminuant la quantité dans le Local Storage"

# Row 9
$ws.Range("B9").Value = "Affichage du prix de l'article"
$ws.Range("C9").Value = "Récupération  du prix, ajouté dans le Local Storage"
$ws.Range("D9").Value = "Le prix de l'article est affiché"

# Row 10
$ws.Range("C10").Value = "Le prix total equivaut au à la somme de chaque articles multiplié par leur quantité"
$ws.Range("D10").Value = "Le total est affiché par l'addition de chaque article"

# Row 11
$ws.Range("C11").Value = "Bouton cliquable supprimant le produit du Local Storage"

# Row 12
$ws.Range("C12").Value = "REGEX, avec messages d'erreurs pour chaque champs"

# Row 13
$ws.Range("C13").Value = "Bouton confirmant la validité du formulaire, et qui effectue une redirection de page"

# Row 14
$ws.Range("C14").Value = "Récupération d'un numéro aléatoire de confirmation de commande du fichier JSON"
$ws.Range("E14").Value = "OK / Si l'on inscrit l'adresse intenret de la page commande, la commande sera indiquée comme étant validée, sans numéro de commande inscrit"
